$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: this registration's data changed (different student / class / recording),
# and the recording link (L3) could not be resolved -> cleared to empty.
$ws.Range("A3").Value = 84410
$ws.Range("B3").Value = 53068043
$ws.Range("C3").Value = "Idiomas: B2:Inglés intermedio, N14 (J-Ad13+)"
$ws.Range("D3").Value = "Daniel Guarin Acevedo"
$ws.Range("E3").Value = "DanielGuarin@comfama.com.co"
$ws.Range("F3").Value = "Clase 20 - 53068043"
$ws.Range("I3").Value = 8300
$ws.Range("J3").Value = "Clase 20 - 53068043"
$ws.Range("K3").Value = "https://comfama.webex.com/comfama/ldr.php?RCID=24172b94001615a9f9924cc1be32e5eddd"
$ws.Range("L3").Value = ""

# Row 4: replaced with the (duplicated) data from row 2's registration.
$ws.Range("A4").Value = 825
$ws.Range("B4").Value = 53068027
$ws.Range("C4").Value = "Idiomas: B2:Inglés intermedio, N12 (J-Ad13+)"
$ws.Range("D4").Value = "EDWIN FERNANDO CARDONA BOCANEGRA"
$ws.Range("E4").Value = "EdwinCardona@comfama.com.co"
$ws.Range("F4").Value = "Clase 19 - 53068027"
$ws.Range("H4").Value = "09:30:00"
$ws.Range("I4").Value = 8215
$ws.Range("J4").Value = "Clase 19 - 53068027"
$ws.Range("K4").Value = "https://comfama.webex.com/comfama/ldr.php?RCID=ca7c616fa8a60ac9d71dbc710d76906c"
$ws.Range("L4").Value = "https://sagrabacionescursos.blob.core.windows.net/videos/Clase 19 - 53068027-20211127 1455-1.mp4"

# Rows 5 and 6 are no longer part of the result set.
$ws.Rows("5:6").Delete()
